$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New exercise rows to append after the existing data (row 34 is last used row)
$newRows = @(
    @{ Name = "subir escadas"; Cal = 10.3 },
    @{ Name = "bicicleta ergométrica"; Cal = 8.3 },
    @{ Name = "dança"; Cal = 6.7 }
)

$startRow = 35
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Carry over the same cell formatting/style used by the existing data rows
    $ws.Range("A34:H34").Copy() | Out-Null
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row.Name
    $ws.Cells.Item($r, 2).Value = "exercícios"
    $ws.Cells.Item($r, 3).Value = "cardio"
    $ws.Cells.Item($r, 4).Value = $row.Cal
    $ws.Cells.Item($r, 5).Value = 0.0
    $ws.Cells.Item($r, 6).Value = 0.0
    $ws.Cells.Item($r, 7).Value = 0.0
    $ws.Cells.Item($r, 8).Value = "exercício"
}

$excel.CutCopyMode = 0
